# Fruta / hortaliza, semanal
#
# Insert one new weekly record for Piña (Vega Monumental Concepción) at
# row 140 — this pushes the existing rows 140-204 down to 141-205 (the
# used range grows from A1:T204 to A1:T205) and fills the freshly
# inserted row with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 140..204 down to 141..205, opening up a blank row 140.
$ws.Rows.Item(140).Insert()

# Populate the new row 140 with the new weekly record.
$ws.Range("A140").Value = 11
$ws.Range("B140").Value = "Vega Monumental Concepción"
$ws.Range("C140").Value = "Bíobío"
$ws.Range("D140").Value = 44825
$ws.Range("E140").Value = 8
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100108
$ws.Range("H140").Value = "Tropicales y subtropicales"
$ws.Range("I140").Value = 100108005
$ws.Range("J140").Value = "Piña"
$ws.Range("K140").Value = "Caramelo"
$ws.Range("L140").Value = "Segunda"
$ws.Range("M140").Value = 140
$ws.Range("N140").Value = 18000
$ws.Range("O140").Value = 20000
$ws.Range("P140").Value = 18857
$ws.Range("Q140").Value = "$/caja 14 unidades"
$ws.Range("R140").Value = "Ecuador"
$ws.Range("S140").Value = 1347
$ws.Range("T140").Value = 14

# Match the date-number formatting used by the rest of column D.
$ws.Range("D140").NumberFormat = $ws.Range("D141").NumberFormat()
